# Add a new "2022-Q1" worksheet (fund holdings for that quarter) between the
# existing "2021-Q4" sheet and the "总计" (summary) sheet, and update the
# "总计" sheet with a new row for 2022-Q1 (pushing the 2021-Q4 row down).

$wb = $excel.ActiveWorkbook

$sheetQ4 = $wb.Worksheets.Item("2021-Q4")

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet right after "2021-Q4"
# ---------------------------------------------------------------------------
$sheetQ1 = $wb.Worksheets.Add($null, $sheetQ4)
$sheetQ1.Name = "2022-Q1"

# NOTE: sheet references behave like position-based pointers, so the "总计"
# sheet must be (re)looked-up *after* the insertion above shifted its index,
# otherwise this variable would keep pointing at whatever sheet now sits in
# its old slot (the newly inserted one).
$sheetSum = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------------
$sheetQ1.Range("B1").Value = "基金代码"
$sheetQ1.Range("C1").Value = "基金名称"
$sheetQ1.Range("D1").Value = "基金规模"
$sheetQ1.Range("E1").Value = "股票总仓位"
$sheetQ1.Range("F1").Value = "仓位占比"
$sheetQ1.Range("G1").Value = "持有市值(亿元)"
$sheetQ1.Range("H1").Value = "仓位排名"

# ---------------------------------------------------------------------------
# 3. Fund holding rows (columns B-G are stored as text, matching the source
#    data which keeps fund codes/percentages as plain strings)
# ---------------------------------------------------------------------------
$fundRows = @(
    @{ Row=2;  Idx=0; Code="011488"; Name="申万菱信乐享混合";                         Size="12.18"; Pos="77.88"; Pct="3.69"; Mv="0.4494"; Rank=8 },
    @{ Row=3;  Idx=1; Code="007202"; Name="天弘优质成长企业精选混合";                 Size="4.81";  Pos="92.52"; Pct="4.49"; Mv="0.2160"; Rank=8 },
    @{ Row=4;  Idx=2; Code="012051"; Name="申万菱信乐道三年持有期混合型证券投资基金"; Size="3.38";  Pos="81.64"; Pct="4.71"; Mv="0.1592"; Rank=5 },
    @{ Row=5;  Idx=3; Code="011851"; Name="天弘先进制造混合型证券投资基金A";          Size="2.72";  Pos="91.41"; Pct="5.73"; Mv="0.1559"; Rank=4 },
    @{ Row=6;  Idx=4; Code="011852"; Name="天弘先进制造混合型证券投资基金C";          Size="0.70";  Pos="91.41"; Pct="5.73"; Mv="0.0401"; Rank=4 },
    @{ Row=7;  Idx=5; Code="004694"; Name="天弘策略精选灵活配置混合A";                Size="1.11";  Pos="80.93"; Pct="3.54"; Mv="0.0393"; Rank=8 },
    @{ Row=8;  Idx=6; Code="000433"; Name="安信鑫发优选混合";                         Size="1.23";  Pos="67.20"; Pct="2.21"; Mv="0.0272"; Rank=7 },
    @{ Row=9;  Idx=7; Code="004748"; Name="天弘策略精选灵活配置混合C";                Size="0.08";  Pos="80.93"; Pct="3.54"; Mv="0.0028"; Rank=8 },
    @{ Row=10; Idx=8; Code="001744"; Name="诺安进取回报灵活配置混合";                 Size="0.04";  Pos="62.10"; Pct="4.38"; Mv="0.0018"; Rank=7 }
)

foreach ($r in $fundRows) {
    $rowNum = $r.Row

    # Keep the text-looking numeric columns (B,D,E,F,G) as real text so that
    # leading zeros / trailing zeros in the source data are preserved.
    $sheetQ1.Range("B$rowNum`:G$rowNum").NumberFormat = "@"

    $sheetQ1.Range("A$rowNum").Value = $r.Idx
    $sheetQ1.Range("B$rowNum").Value = $r.Code
    $sheetQ1.Range("C$rowNum").Value = $r.Name
    $sheetQ1.Range("D$rowNum").Value = $r.Size
    $sheetQ1.Range("E$rowNum").Value = $r.Pos
    $sheetQ1.Range("F$rowNum").Value = $r.Pct
    $sheetQ1.Range("G$rowNum").Value = $r.Mv
    $sheetQ1.Range("H$rowNum").Value = $r.Rank
}

# ---------------------------------------------------------------------------
# 4. Update the "总计" (summary) sheet: push the existing 2021-Q4 row down to
#    row 3 (copying it so formatting/style is preserved) and insert a new
#    2022-Q1 row in its place at row 2.
# ---------------------------------------------------------------------------
$sheetSum.Range("A2:D2").Copy($sheetSum.Range("A3:D3"))
$sheetSum.Range("A3").Value = 1

$sheetSum.Range("A2").Value = 0
$sheetSum.Range("B2").Value = "2022-Q1"
$sheetSum.Range("C2").Value = 9
$sheetSum.Range("D2").Value = 1.09

Write-Host "2022-Q1 sheet added and 总计 sheet updated"
